# TAC-3804  Fix allow carrier saas to import trip , edit excel files
#
# The sample import template had a stray "Dry goods" default value sitting
# in C2 (under the "Goods Sub Category" header) and the Unit Of Measure
# dropdown for column G only offered a reduced set of choices. This cleans
# that up:
#   - removes the leftover sample value from C2 on Sheet1
#   - replaces the G2:G1048576 list validation with the fuller option set
#   - leaves the selection sitting on G1 (the Unit Of Measure header)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stray "Dry goods" sample value left in C2 (second row of the
# template). Nothing else lives on row 2, so the row disappears entirely
# and the sheet's used range shrinks back down to A1:M1.
$ws.Range("C2").ClearContents()

# Replace the Unit Of Measure (column G) dropdown list with the expanded
# set of choices.
$ws.Range("G2:G1048576").Validation.Delete()
$ws.Range("G2:G1048576").Validation.Add(3, 1, 1, '"Litre,Box,Bag,Piece,Weight -KG,Pallets,Container,Others"')

# Move the active selection / view over to the Unit Of Measure column.
$ws.Range("G1").Select()
